$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "published" sheet: the tuberculosis-screening paper (previously in
#    "submitted") has now been published. Append it as a new row to the
#    Tabel1 table, copying the formatting used for earlier "just added"
#    rows (e.g. row 40) so the highlighted style (s=14/15/16/17/18)
#    carries over correctly.
# ---------------------------------------------------------------------
$wsPub = $wb.Worksheets.Item("published")
$loPub = $wsPub.ListObjects.Item(1)
$loPub.ListRows.Add() | Out-Null

$srcRow = $wsPub.Range("A40:M40")
$dstRow = $wsPub.Range("A52:M52")
$srcRow.Copy($dstRow)

$wsPub.Range("A52").Value = "Cost-effectiveness of tuberculosis screening policies in Flanders, Belgium"
$wsPub.Range("B52").Value = "Smit, G Suzanne A; Apers, Ludwig; Arrazole de Onate, Wouter; Beutels, Philippe; Dorny, Pierre; Forier, An-Marie; Janssens, Kristien; Macq, Jean; Mak, Ruud; Schol, Sandrina; Wildemeersch, Dirk; Speybroeck, Niko; Devleesschauwer, Brecht"
$wsPub.Range("C52").Value = "Bulletin of the World Health Organization"
$wsPub.Range("D52").Value = "Bull. World Health Org."
$wsPub.Range("E52").Value = 2016
$wsPub.Range("F52").Value = "NA"
$wsPub.Range("G52").Value = "NA"
$wsPub.Range("H52").Value = "NA"
$wsPub.Range("I52").Value = "NA"
$wsPub.Range("J52").Value = "NA"
$wsPub.Range("K52").Value = 42641
$wsPub.Range("L52").Value = ""
$wsPub.Range("M52").Value = ""

# ---------------------------------------------------------------------
# 2. "submitted" sheet: remove the now-published tuberculosis-screening
#    row (row 5) - everything below shifts up automatically - then add
#    a brand-new submission (the mercury-intoxication disability-weights
#    paper) as the new last row.
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("submitted")
$wsSub.Rows.Item(5).Delete()
$loSub = $wsSub.ListObjects.Item(1)
$loSub.ListRows.Add() | Out-Null

$srcRow2 = $wsSub.Range("A13:L13")
$dstRow2 = $wsSub.Range("A14:L14")
$srcRow2.Copy($dstRow2)

$wsSub.Range("A14").Value = "Disability weights for chronic mercury intoxication resulting from gold mining activities: results from an online pairwise comparisons survey"
$wsSub.Range("C14").Value = "Environmental Health Perspectives"
$wsSub.Range("D14").Value = "Environ. Health Perspect."
$wsSub.Range("B14").Value = "Steckling, Nadine; Devleesschauwer, Brecht; Winkelnkemper, Julia; Fischer, Florian; Ericson, Bret; Krämer, Alexander; Hornberg, Claudia; Fuller, Richard; Plass, Dietrich; Bose-O'Reilly, Stephan"
$wsSub.Range("K14").Value = 42644
$wsSub.Range("L14").Formula = "=TODAY()-K14"
